$wb = $excel.ActiveWorkbook
$wsPonds = $wb.Worksheets.Item("Ponds")
$wsEggrooms = $wb.Worksheets.Item("Eggrooms")

# Remove the "Stock 1", "Stock 2", "Stock 3" columns (E:G) on both sheets -
# the parsers no longer read per-stock columns now that distributions /
# treatments are parsed by dedicated classes.
$wsPonds.Range("E1:G1").EntireColumn.Delete()
$wsEggrooms.Range("E1:G1").EntireColumn.Delete()

# Ponds sheet gained a tall blank header row above the table.
$wsPonds.Rows("1").RowHeight = 51
$wsPonds.Rows("3").RowHeight = 45.75

# Update selections on each sheet.
$wsPonds.Range("J4").Select()
$wsEggrooms.Range("G6").Select()
